$d = $word.ActiveDocument

# --- Paragraph: "二维表新增字段调整不了字段的信息" -------------------------
# Split into two runs around the re-inserted _GoBack bookmark and
# highlight both halves red.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("二维表新增字段调整不了字段", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
  $rng.HighlightColorIndex = 6
  $bm = $rng.Duplicate
  $bm.Collapse(0)
  $d.Bookmarks.Add("_GoBack", $bm)
}

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("的信息", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
  $rng.HighlightColorIndex = 6
}

# --- Highlight a handful of the remaining todo-list bullets red -----------
$toHighlight = @(
  "二维表热更新增字段的时候，保存不了顺序",
  "单元格可以用enter输入换行",
  "展开界面，字段能在nil的时候保存",
  "筛选的时候，多个字段筛选没有根据之前的筛选结果二次筛选"
)

foreach ($text in $toHighlight) {
  $rng = $d.Content
  $rng.Find.ClearFormatting()
  $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
  if ($found) {
    $rng.HighlightColorIndex = 6
  }
}
